$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new header columns: W (finalVerdict) and X (finalVerdictDate) ---

$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Range("V1").Copy()
$ws.Range("X1").PasteSpecial(-4122)

$ws.Range("W1").Value = "الحكم النهائي "

$ws.Range("X1").Value = "تاريخ الحكم النهائي (dd/mm/yyyy) "
$run1 = $ws.Range("X1").Characters(1, 20)
$run1.Font.Bold = $true
$run1.Font.Size = 14
$run2 = $ws.Range("X1").Characters(21)
$run2.Font.Bold = $true
$run2.Font.Size = 10

# --- Simplify old per-column date formats down to plain text (@) ---
$ws.Range("H1").NumberFormat = "@"
$ws.Range("K1").NumberFormat = "@"
$ws.Range("O1").NumberFormat = "@"
$ws.Range("R1").NumberFormat = "@"

# --- Row height tweaks (rows 218/219 grow from 13.5 to 15.75) ---
$ws.Rows.Item(218).RowHeight = 15.75
$ws.Rows.Item(219).RowHeight = 15.75
